$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each updated cell to remain text (matches original inlineStr string values)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.59%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.241"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.96%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07552"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.938"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.27%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.818"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.30%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.526"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "8.97%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9234"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.39%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1693"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07879"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.76%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08023"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.73%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03015"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.43%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09913"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "10.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-5.44%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.40%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006398"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.59%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.447"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.11%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.235"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.07%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.46%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1327"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.02%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.474"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.66%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1620"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.81%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.93%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004457"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.75%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001399"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "20.00%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01692"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2,484.40%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04471"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.68%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006971"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.44%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1353"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.78%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002078"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01373"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.02%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006157"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.46%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7090"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-63.25%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.44%"
